$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.529.49"
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = "'1.736.77"
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'246.68"
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = "'0.4929"
$ws.Range("E7").Value = '  +2.91%  '
$ws.Range("D8").Value = "'0.2665"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = "'0.06286"
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").Value = "'1.732.95"
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").Value = "'0.07044"
$ws.Range("E11").Value = '  -1.17%  '
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = "'4.597"
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = "'0.6114"
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").Value = "'77.47"
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = "'0.000007352"
$ws.Range("E17").Value = '  +6.54%  '
$ws.Range("D18").Value = "'26.531.83"
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  -1.39%  '
$ws.Range("D21").Value = "'1.956.15"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = "'4.590"
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").Value = "'8.700"
$ws.Range("E23").Value = '  -2.20%  '
$ws.Range("D24").Value = "'5.250"
$ws.Range("D25").Value = "'139.89"
$ws.Range("E25").Value = '  +2.67%  '
$ws.Range("D26").Value = "'15.45"
$ws.Range("E26").Value = '  +0.69%  '
$ws.Range("D27").Value = "'1.420"
$ws.Range("E27").Value = '  +1.48%  '
$ws.Range("D28").Value = "'108.05"
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("D29").Value = "'1.762"
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").Value = "'4.038"
$ws.Range("E30").Value = '  +1.51%  '
$ws.Range("D31").Value = "'0.08064"
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("D32").Value = "'3.717"
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("D33").Value = "'0.04597"
$ws.Range("E33").Value = '  +1.33%  '
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").Value = "'1.008"
$ws.Range("E35").Value = '  +1.92%  '
$ws.Range("D36").Value = "'0.6363"
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").Value = "'0.8971"
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("D38").Value = "'2.013"
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("D39").Value = "'2.403"
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("D41").Value = "'0.01509"
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("D42").Value = "'101.92"
$ws.Range("E42").Value = '  -7.07%  '
$ws.Range("D43").Value = "'5.402"
$ws.Range("E43").Value = '  -4.40%  '
$ws.Range("D44").Value = "'0.3903"
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("D45").Value = "'6.865"
$ws.Range("E45").Value = '  -1.04%  '
$ws.Range("D46").Value = "'0.1187"
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("E47").Value = '  +1.31%  '
$ws.Range("D48").Value = "'30.54"
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("D49").Value = "'7.777"
$ws.Range("E49").Value = '  -1.13%  '
$ws.Range("D50").Value = "'1.268"
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").Value = "'51.78"
